$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text of A5: it previously described "(Half-rate F32 Accumulate
# during training)" and is replaced with a more detailed note about which
# precisions are affected.
$ws.Range("A5").Value = "(Half-rate FP16/BF16/TF32 if using FP32 Accumulate during training)"

# Add a new trailing note in A23 clarifying the author doesn't own a RTX 3090.
$ws.Range("A23").Value = "Note: I do not have a RTX 3090."

# Move the active selection to reflect where the author ended up editing.
$ws.Range("A24").Select()
